$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '20.605.87'
$ws.Range('E2').Value = '  +2.53%  '

# Row 3
$ws.Range('D3').Value = '1.476.27'
$ws.Range('E3').Value = '  +3.29%  '

# Row 4
$ws.Range('E4').Value = '  +0.63%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9893'
$ws.Range('E5').Value = '  -1.19%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '281.10'
$ws.Range('E6').Value = '  +2.23%  '

# Row 7
$ws.Range('E7').Value = '  +1.10%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3222'
$ws.Range('E8').Value = '  +4.78%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '41.89'
$ws.Range('E9').Value = '  +4.44%  '

# Row 10
$ws.Range('E10').Value = '  +7.10%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06763'
$ws.Range('E11').Value = '  +2.98%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  -0.06%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.700'
$ws.Range('E13').Value = '  +5.10%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.65'
$ws.Range('E14').Value = '  +8.23%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.352'
$ws.Range('E15').Value = '  +2.88%  '

# Row 16
$ws.Range('D16').Value = '1.475.99'
$ws.Range('E16').Value = '  +2.54%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001045'
$ws.Range('E17').Value = '  +3.33%  '

# Row 18
$ws.Range('E18').Value = '  +0.10%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.48'
$ws.Range('E19').Value = '  -2.56%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9869'
$ws.Range('E20').Value = '  -1.44%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.754'
$ws.Range('E21').Value = '  +1.26%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.03'
$ws.Range('E22').Value = '  +3.72%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.31'
$ws.Range('E23').Value = '  +1.95%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.315'
$ws.Range('E24').Value = '  -0.92%  '

# Row 25
$ws.Range('D25').Value = '20.784.41'
$ws.Range('E25').Value = '  +3.33%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.352'
$ws.Range('E26').Value = '  +2.76%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '138.92'
$ws.Range('E27').Value = '  +0.25%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.82'
$ws.Range('E28').Value = '  +5.46%  '

# Row 29
$ws.Range('D29').Value = '1.641.43'
$ws.Range('E29').Value = '  +2.95%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.08'
$ws.Range('E30').Value = '  +5.12%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.002'
$ws.Range('E31').Value = '  +2.50%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.434'
$ws.Range('E32').Value = '  -0.09%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8566'
$ws.Range('E33').Value = '  -6.83%  '

# Row 34: 'Stellar' -> 'WEMIXTOKEN'
$ws.Range('B34').Value = 'WEMIXTOKEN'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.632'
$ws.Range('E34').Value = '  +25.03%  '

# Row 35: 'WEMIXTOKEN' -> 'Stellar'
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07904'
$ws.Range('E35').Value = '  +1.68%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06110'
$ws.Range('E36').Value = '  +7.60%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.003'
$ws.Range('E37').Value = '  +4.41%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '10.88'
$ws.Range('E38').Value = '  -4.77%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9923'
$ws.Range('E39').Value = '  -0.87%  '

# Row 40: 'VeChain' -> 'FraxShare'
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.866'
$ws.Range('E40').Value = '  -6.80%  '

# Row 41: 'FraxShare' -> 'VeChain'
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.02094'
$ws.Range('E41').Value = '  +4.02%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.143'
$ws.Range('E42').Value = '  +2.29%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1925'
$ws.Range('E43').Value = '  +0.36%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5485'
$ws.Range('E44').Value = '  +3.11%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.58'
$ws.Range('E45').Value = '  +2.11%  '

# Row 46
$ws.Range('E46').Value = '  +1.45%  '

# Row 47: 'Decentraland' -> 'Quant'
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '122.42'
$ws.Range('E47').Value = '  +10.48%  '

# Row 48: 'Quant' -> 'Decentraland'
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5422'
$ws.Range('E48').Value = '  +5.75%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.846'
$ws.Range('E49').Value = '  +3.56%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.067'
$ws.Range('E50').Value = '  +1.50%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06483'
$ws.Range('E51').Value = '  +4.42%  '
